$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: for numeric-looking text values in column D, force text type
# via NumberFormat="@" then restore the cell style to Normal so no stray
# number-format style is left attached to the cell (matches original OOXML).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.999.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -6.49%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.293.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.38%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.86%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.294.48"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.33%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.470"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.74%  "

$ws.Range("E11").Value = "  -5.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.370"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.30%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.857.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.26%  "

$ws.Range("E14").Value = "  -0.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.297.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.19%  "

$ws.Range("E16").Value = "  -6.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.217.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "24.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.45%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -11.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "348.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -9.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.552"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.67%  "

$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.425.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "68.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.87%  "

$ws.Range("E27").Value = "  -4.37%  "

$ws.Range("E28").Value = "  +0.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.151"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.00%  "

$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.323.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.93%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.72"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.80%  "

$ws.Range("E39").Value = "  -2.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "156.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0747"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.19%  "

$ws.Range("E42").Value = "  +0.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.739"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.48%  "

$ws.Range("E46").Value = "  +2.38%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.72%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.71%  "

$ws.Range("E49").Value = "  -0.96%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.855"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.18%  "
